$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1987767584097859
$ws.Range("C2").Value = 0.5504587155963303
$ws.Range("J2").Value = 0.01529051987767584
$ws.Range("O2").Value = 0.009174311926605505
$ws.Range("P2").Value = 0.1284403669724771
$ws.Range("S2").Value = 0.09785932721712538
$ws.Range("B3").Value = 0.02150537634408602
$ws.Range("C3").Value = 0.02688172043010753
$ws.Range("J3").Value = 0.01612903225806452
$ws.Range("P3").Value = 0.7204301075268817
$ws.Range("S3").Value = 0.2150537634408602
$ws.Range("J4").Value = 0.09259259259259259
$ws.Range("P4").Value = 0.5555555555555556
$ws.Range("S4").Value = 0.3518518518518519
$ws.Range("B6").Value = 0.05371900826446281
$ws.Range("D6").Value = 0.02066115702479339
$ws.Range("F6").Value = 0.04545454545454546
$ws.Range("J6").Value = 0.2355371900826446
$ws.Range("O6").Value = 0.04545454545454546
$ws.Range("Q6").Value = 0.2024793388429752
$ws.Range("R6").Value = 0.05785123966942149
$ws.Range("S6").Value = 0.3388429752066116
$ws.Range("B7").Value = 0.1280788177339902
$ws.Range("D7").Value = 0.03448275862068965
$ws.Range("F7").Value = 0.0541871921182266
$ws.Range("J7").Value = 0.1921182266009852
$ws.Range("O7").Value = 0.004926108374384237
$ws.Range("Q7").Value = 0.1625615763546798
$ws.Range("R7").Value = 0.07389162561576355
$ws.Range("S7").Value = 0.3497536945812808
$ws.Range("B8").Value = 0.09003831417624521
$ws.Range("D8").Value = 0.02298850574712644
$ws.Range("E8").Value = 0.001915708812260536
$ws.Range("F8").Value = 0.06513409961685823
$ws.Range("J8").Value = 0.1130268199233716
$ws.Range("O8").Value = 0.01724137931034483
$ws.Range("Q8").Value = 0.1609195402298851
$ws.Range("R8").Value = 0.06704980842911877
$ws.Range("S8").Value = 0.4616858237547893
$ws.Range("B10").Value = 0.1155261221486387
$ws.Range("D10").Value = 0.02281089036055924
$ws.Range("E10").Value = 0.0007358351729212656
$ws.Range("F10").Value = 0.07652685798381163
$ws.Range("J10").Value = 0.1037527593818985
$ws.Range("O10").Value = 0.02060338484179544
$ws.Range("Q10").Value = 0.1905813097866078
$ws.Range("R10").Value = 0.07799852832965416
$ws.Range("S10").Value = 0.3914643119941133
$ws.Range("F11").Value = 0.002857142857142857
$ws.Range("G11").Value = 0.1142857142857143
$ws.Range("J11").Value = 0.1485714285714286
$ws.Range("K11").Value = 0.1885714285714286
$ws.Range("L11").Value = 0.5371428571428571
$ws.Range("S11").Value = 0.008571428571428572
$ws.Range("G12").Value = 0.6994818652849741
$ws.Range("J12").Value = 0.227979274611399
$ws.Range("K12").Value = 0.02072538860103627
$ws.Range("L12").Value = 0.0310880829015544
$ws.Range("S12").Value = 0.02072538860103627
$ws.Range("G13").Value = 0.6326530612244898
$ws.Range("J13").Value = 0.3673469387755102
$ws.Range("F15").Value = 0.01234567901234568
$ws.Range("H15").Value = 0.1810699588477366
$ws.Range("I15").Value = 0.06172839506172839
$ws.Range("J15").Value = 0.3415637860082305
$ws.Range("K15").Value = 0.06172839506172839
$ws.Range("M15").Value = 0.02880658436213992
$ws.Range("O15").Value = 0.08230452674897119
$ws.Range("S15").Value = 0.2304526748971193
$ws.Range("F16").Value = 0.03
$ws.Range("H16").Value = 0.125
$ws.Range("I16").Value = 0.06
$ws.Range("J16").Value = 0.45
$ws.Range("K16").Value = 0.12
$ws.Range("M16").Value = 0.02
$ws.Range("O16").Value = 0.05
$ws.Range("S16").Value = 0.145
$ws.Range("F17").Value = 0.0244988864142539
$ws.Range("H17").Value = 0.22271714922049
$ws.Range("I17").Value = 0.04899777282850779
$ws.Range("J17").Value = 0.4008908685968819
$ws.Range("K17").Value = 0.111358574610245
$ws.Range("M17").Value = 0.0155902004454343
$ws.Range("O17").Value = 0.0556792873051225
$ws.Range("S17").Value = 0.1202672605790646
$ws.Range("F18").Value = 0.03932584269662921
$ws.Range("H18").Value = 0.1741573033707865
$ws.Range("I18").Value = 0.07865168539325842
$ws.Range("J18").Value = 0.4269662921348314
$ws.Range("K18").Value = 0.0898876404494382
$ws.Range("M18").Value = 0.01685393258426966
$ws.Range("O18").Value = 0.06179775280898876
$ws.Range("S18").Value = 0.1123595505617977
$ws.Range("F19").Value = 0.01204819277108434
$ws.Range("H19").Value = 0.2303330970942594
$ws.Range("I19").Value = 0.05315379163713678
$ws.Range("J19").Value = 0.367115520907158
$ws.Range("K19").Value = 0.1226080793763288
$ws.Range("M19").Value = 0.02055279943302622
$ws.Range("N19").Value = 0.001417434443656981
$ws.Range("O19").Value = 0.06732813607370659
$ws.Range("S19").Value = 0.1254429482636428
